# --------------------------------------------------------------------------
# OLX Monitor 2026-02-28 07:48
# Append 9 freshly-scraped OLX listings (rows 252-260) to the "PODSUMOWANIE"
# sheet, matching the layout/formatting of the existing monitoring rows.
# --------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PODSUMOWANIE")

# Template row supplies the base formatting for every new row (cols A-H).
$template = $ws.Range("A251:H251")

# Donor cells used to stamp the two non-default "days" (col F) styles that show
# up below, without ever creating/altering entries in the styles table:
$fStyleNew = $ws.Range("F234")   # style for freshly (re)posted listings, s=16
$fStyleMid = $ws.Range("F7")     # style for medium-age listings, s=14
# (the template row above already carries the "old listing" style, s=15,
#  which is what most of the rows below need)

# NOTE: column E holds dates formatted as free text (e.g. "10.10.2025"). Typing
# those straight into .Value would get silently re-interpreted by Excel as a
# real date (and reformatted/serialized) whenever day <= 12. To guarantee they
# stay literal text - exactly like the rest of the sheet - each one is copied
# from another cell in the sheet that already holds that exact string as text.

# --- Row 252 ---
$template.Copy($ws.Range("A252:H252"))
$ws.Range("A252").Value = '2026-02-28 07:48:03'
$ws.Range("B252").Value = 'poqui'
$ws.Range("C252").Value = 'Kawalerka po remoncie z funkcjonalną antresolą - ul. Jana Sawy'
$ws.Range("D252").Value = 2499
$ws.Range("E8").Copy($ws.Range("E252"))   # "28.10.2025" as text
$ws.Range("F252").Value = 122
$ws.Range("G252").Value = 'https://www.olx.pl/d/oferta/kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger.html'
$ws.Range("H252").Value = 'kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger'

# --- Row 253 ---
$template.Copy($ws.Range("A253:H253"))
$ws.Range("A253").Value = '2026-02-28 07:48:03'
$ws.Range("B253").Value = 'poqui'
$ws.Range("C253").Value = 'Duży pokój z balkonem w 2pokojowym mieszkaniu blisko Politechniki'
$ws.Range("D253").Value = 1665
$ws.Range("E234").Copy($ws.Range("E253"))   # "25.02.2026" as text
$fStyleNew.Copy()
$ws.Range("F253").PasteSpecial(-4122)
$ws.Range("F253").Value = 2
$ws.Range("G253").Value = 'https://www.olx.pl/d/oferta/duzy-pokoj-z-balkonem-w-2pokojowym-mieszkaniu-blisko-politechniki-CID3-ID19xpQK.html'
$ws.Range("H253").Value = 'duzy-pokoj-z-balkonem-w-2pokojowym-mieszkaniu-blisko-politechniki-CID3-ID19xpQK'

# --- Row 254 ---
$template.Copy($ws.Range("A254:H254"))
$ws.Range("A254").Value = '2026-02-28 07:48:03'
$ws.Range("B254").Value = 'poqui'
$ws.Range("C254").Value = 'Nowoczesne mieszkanie 2-pokojowe z balkonem, blisko UMCS, KUL, UP'
$ws.Range("D254").Value = 2499
$ws.Range("E234").Copy($ws.Range("E254"))   # "25.02.2026" as text
$fStyleNew.Copy()
$ws.Range("F254").PasteSpecial(-4122)
$ws.Range("F254").Value = 2
$ws.Range("G254").Value = 'https://www.olx.pl/d/oferta/nowoczesne-mieszkanie-2-pokojowe-z-balkonem-blisko-umcs-kul-up-CID3-ID19xpwN.html'
$ws.Range("H254").Value = 'nowoczesne-mieszkanie-2-pokojowe-z-balkonem-blisko-umcs-kul-up-CID3-ID19xpwN'

# --- Row 255 ---
$template.Copy($ws.Range("A255:H255"))
$ws.Range("A255").Value = '2026-02-28 07:48:03'
$ws.Range("B255").Value = 'poqui'
$ws.Range("C255").Value = 'Przytulny pokój blisko Politechniki – ul. Przytulna'
$ws.Range("D255").Value = 549
$ws.Range("E9").Copy($ws.Range("E255"))   # "10.10.2025" as text
$ws.Range("F255").Value = 140
$ws.Range("G255").Value = 'https://www.olx.pl/d/oferta/przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz.html'
$ws.Range("H255").Value = 'przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz'

# --- Row 256 ---
$template.Copy($ws.Range("A256:H256"))
$ws.Range("A256").Value = '2026-02-28 07:48:03'
$ws.Range("B256").Value = 'poqui'
$ws.Range("C256").Value = 'Mieszkanie z KLIMATYZACJĄ 5 minut od UMCS, UP, KUL - Długosza'
$ws.Range("D256").Value = 2049
$ws.Range("E10").Copy($ws.Range("E256"))   # "19.12.2025" as text
$ws.Range("F256").Value = 70
$ws.Range("G256").Value = 'https://www.olx.pl/d/oferta/mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc.html'
$ws.Range("H256").Value = 'mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc'

# --- Row 257 ---
$template.Copy($ws.Range("A257:H257"))
$ws.Range("A257").Value = '2026-02-28 07:48:03'
$ws.Range("B257").Value = 'pokojewlublinie'
$ws.Range("C257").Value = 'WOLNY OD ZARAZ! Super lokalizacja, blisko centrum, ul. Paganiniego 12'
$ws.Range("D257").Value = 12640
$ws.Range("E7").Copy($ws.Range("E257"))   # "19.01.2026" as text
$fStyleMid.Copy()
$ws.Range("F257").PasteSpecial(-4122)
$ws.Range("F257").Value = 39
$ws.Range("G257").Value = 'https://www.olx.pl/d/oferta/wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc.html'
$ws.Range("H257").Value = 'wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc'

# --- Row 258 ---
$template.Copy($ws.Range("A258:H258"))
$ws.Range("A258").Value = '2026-02-28 07:48:03'
$ws.Range("B258").Value = 'pokojewlublinie'
$ws.Range("C258").Value = 'WOLNY OD ZARAZ! Pokój jedynka, ul. Romanowskiego 58'
$ws.Range("D258").Value = 0
$ws.Range("E11").Copy($ws.Range("E258"))   # "11.08.2025" as text
$ws.Range("F258").Value = 200
$ws.Range("G258").Value = 'https://www.olx.pl/d/oferta/wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm.html'
$ws.Range("H258").Value = 'wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm'

# --- Row 259 ---
$template.Copy($ws.Range("A259:H259"))
$ws.Range("A259").Value = '2026-02-28 07:48:03'
$ws.Range("B259").Value = 'dawnypatron'
$ws.Range("C259").Value = 'Ładny pokój jednoosobowy. Wynajmę duży pokój w centrum. ul Niecała 4.'
$ws.Range("D259").Value = 730
$ws.Range("E13").Copy($ws.Range("E259"))   # "20.09.2024" as text
$ws.Range("F259").Value = 525
$ws.Range("G259").Value = 'https://www.olx.pl/d/oferta/ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM.html'
$ws.Range("H259").Value = 'ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM'

# --- Row 260 ---
$template.Copy($ws.Range("A260:H260"))
$ws.Range("A260").Value = '2026-02-28 07:48:03'
$ws.Range("B260").Value = 'dawnypatron'
$ws.Range("C260").Value = 'Mam do wynajęcia pokój dla os. pracującej lub studenta. Narutowicza 14'
$ws.Range("D260").Value = 14690
$ws.Range("E14").Copy($ws.Range("E260"))   # "05.12.2025" as text
$ws.Range("F260").Value = 84
$ws.Range("G260").Value = 'https://www.olx.pl/d/oferta/mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv.html'
$ws.Range("H260").Value = 'mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv'

$excel.CutCopyMode = $false

Write-Host "Added rows 252-260 to PODSUMOWANIE sheet"
